$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27: drop the duplicate "Periodicity not specified" row — shift the
# field_frequency "Quarter" row up into row 27.
$ws.Cells.Item(27, 1).Value = "field_frequency"
$ws.Cells.Item(27, 2).Value = "Quarter"
$ws.Cells.Item(27, 3).Value = "Quarter"

# Rows 28-31: field_wbddh_data_type block (shifted up by one row).
$ws.Cells.Item(28, 1).Value = "field_wbddh_data_type"
$ws.Cells.Item(28, 2).Value = "Cross sectional"
$ws.Cells.Item(28, 3).Value = "Other"

$ws.Cells.Item(29, 1).Value = "field_wbddh_data_type"
$ws.Cells.Item(29, 2).Value = "Time Series"
$ws.Cells.Item(29, 3).Value = "Time Series"

$ws.Cells.Item(30, 1).Value = "field_wbddh_data_type"
$ws.Cells.Item(30, 2).Value = "Transactions"
$ws.Cells.Item(30, 3).Value = "Other"

$ws.Cells.Item(31, 1).Value = "field_wbddh_data_type"
$ws.Cells.Item(31, 2).Value = "Survey(Microdata)"
$ws.Cells.Item(31, 3).Value = "Other"

# New field_granularity_list block — rows 32-36.
$ws.Cells.Item(32, 1).Value = "field_granularity_list"
$ws.Cells.Item(32, 2).Value = "Other"
$ws.Cells.Item(32, 3).Value = "Other"

$ws.Cells.Item(33, 1).Value = "field_granularity_list"
$ws.Cells.Item(33, 2).Value = "National"
$ws.Cells.Item(33, 3).Value = "National"

$ws.Cells.Item(34, 1).Value = "field_granularity_list"
$ws.Cells.Item(34, 2).Value = "Project"
$ws.Cells.Item(34, 3).Value = "Project"

$ws.Cells.Item(35, 1).Value = "field_granularity_list"
$ws.Cells.Item(35, 2).Value = "Sub-national"

$ws.Cells.Item(36, 1).Value = "field_granularity_list"
$ws.Cells.Item(36, 2).Value = "Regional"
$ws.Cells.Item(36, 3).Value = "Regional"

# "Subnational" (row 35, col C) is written last so this new shared string is
# registered after "Regional" — matching the source workbook's string order.
$ws.Cells.Item(35, 3).Value = "Subnational"

# Refresh the view to match the saved window state: scrolled so row 19 is
# the top row, with C37 selected just past the new last data row.
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("C37").Select()
